$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "BodyFat (%)" column header in H1 (must be entered before any other
# new string so it lands first in the shared-string table).
$ws.Range("H1").Value = "BodyFat (%)"
$ws.Range("H1").NumberFormat = "0.0"

# Rows 2-9 don't have a BodyFat reading yet -> literal "NULL" placeholder text.
$ws.Range("H2").Value = "NULL"
$ws.Range("H3").Value = "NULL"
$ws.Range("H4").Value = "NULL"
$ws.Range("H5").Value = "NULL"
$ws.Range("H6").Value = "NULL"
$ws.Range("H7").Value = "NULL"
$ws.Range("H8").Value = "NULL"
$ws.Range("H9").Value = "NULL"
$ws.Range("H9").NumberFormat = "0.0"

# Rows 10-16 get real BodyFat % readings.
$ws.Range("H10").Value = 18.6
$ws.Range("H11").Value = 20.8
$ws.Range("H12").Value = 20.8
$ws.Range("H13").Value = 20.8
$ws.Range("H14").Value = 20.6
$ws.Range("H15").Value = 20.5
$ws.Range("H16").Value = 20.6
$ws.Range("H10:H16").NumberFormat = "0.0"

# Match the column width Excel picked when it auto-fit the new column.
$ws.Columns.Item(8).ColumnWidth = 9.9

# Restore the selection/scroll state recorded in the saved view.
$ws.Range("H15").Select()
